# traded, fixed issues with the repeater
# Append the new trade row (row 8) below the existing trade history,
# copying formats from the row above so the date (col A) and the
# IsShortSell flag (col G) keep the same number formatting/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)

$ws.Range("A8").Value = 42654.743726851855
$ws.Range("B8").Value = $false
$ws.Range("C8").Value = 9835.3700000000008
$ws.Range("D8").Value = 9842.75
$ws.Range("E8").Value = 104.43
$ws.Range("F8").Value = 104.269997
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = -0.15
$ws.Range("I8").Value = $false

# Column A widened slightly to fit the new, slightly longer value.
$ws.Columns("A").ColumnWidth = 14.5
